$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.329.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.15%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.521.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.22%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.62%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.31%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.28%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.64%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.71%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.08%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.25%  "

# Row 13
$ws.Range("E13").Value = "  -5.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.895.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.38%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.476.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.68%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.835"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.12%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.305.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.30%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.39%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0952"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.29%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.02%  "

# Row 24
$ws.Range("E24").Value = "  -1.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.92%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.88%  "

# Row 27
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.34%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.96%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.68%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.04%  "

# Row 33
$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.85%  "

# Row 34
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.84%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0788"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.43%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.53%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.111"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.00%  "

# Row 40
$ws.Range("E40").Value = "  -0.53%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.12%  "

# Row 42
$ws.Range("E42").Value = "  -0.24%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.10%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.979.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.81%  "

# Row 47
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "83.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.63%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.749.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.63%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.90%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.16%  "
